# Apply "minor update to expertise": rename the expertise tag
# "Agroforestry farmer" -> "Farmer" wherever it appears in column J
# (the "Expertise" column) across all data worksheets, and move the
# active selection from Sheet1!J6 to Sheet2!F24 (making Sheet2 the
# active/selected tab).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $firstRow = $used.Row
    $lastRow = $firstRow + $rowCount - 1

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 10)   # column J = "Expertise"
        $text = $cell.Text
        if ($text -and $text.Contains("Agroforestry farmer")) {
            $newText = $text.Replace("Agroforestry farmer", "Farmer")
            $cell.Value = $newText
        }
    }
}

# Switch the active tab/selection from Sheet1 to Sheet2, landing on F24
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("F24").Select()
